# Generate Report for Handback
#
# Row 7 in both the "zh-cn" and "de-de" sheets describes the handback for
# 2944a92a-7419-4e35-862d-fd9e9e20fc8d.md. This run discovered a target
# xlf file and a handback datetime for each locale, but the handback was
# stale (not built from the latest source commit), so an error message is
# recorded as well.

$wb = $excel.ActiveWorkbook

# ---- zh-cn sheet (row 7) ------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# I7 = "Latest Target File" -> becomes a hyperlink to the md file, same as
# the existing A7 / I5 style hyperlinks already on the sheet.
$wsZh.Hyperlinks.Add(
    $wsZh.Range("I7"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/0d18a7ef0f372ac36effbfbd95cccab631602457/e2e/2944a92a-7419-4e35-862d-fd9e9e20fc8d.md",
    "",
    "",
    "2944a92a-7419-4e35-862d-fd9e9e20fc8d.md"
)

# J7 = "Latest Handback File"
$wsZh.Range("J7").Value = "2944a92a-7419-4e35-862d-fd9e9e20fc8d.87cbd9f887a88f593d0c3eb305b24a8a1edbe8ce.zh-cn.xlf"
# K7 = "Latest Handback DateTime"
$wsZh.Range("K7").Value = "2016-08-25 06:53:43"
# P7 = "Error Detail"
$wsZh.Range("P7").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8fe988a4805de0368cba065ccabf6b9f55a68e0a/e2e/2944a92a-7419-4e35-862d-fd9e9e20fc8d.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0d18a7ef0f372ac36effbfbd95cccab631602457/e2e/2944a92a-7419-4e35-862d-fd9e9e20fc8d.md."

# ---- de-de sheet (row 7) --------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Add(
    $wsDe.Range("I7"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/0d18a7ef0f372ac36effbfbd95cccab631602457/e2e/2944a92a-7419-4e35-862d-fd9e9e20fc8d.md",
    "",
    "",
    "2944a92a-7419-4e35-862d-fd9e9e20fc8d.md"
)

$wsDe.Range("J7").Value = "2944a92a-7419-4e35-862d-fd9e9e20fc8d.87cbd9f887a88f593d0c3eb305b24a8a1edbe8ce.de-de.xlf"
$wsDe.Range("K7").Value = "2016-08-25 06:53:50"
$wsDe.Range("P7").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8fe988a4805de0368cba065ccabf6b9f55a68e0a/e2e/2944a92a-7419-4e35-862d-fd9e9e20fc8d.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0d18a7ef0f372ac36effbfbd95cccab631602457/e2e/2944a92a-7419-4e35-862d-fd9e9e20fc8d.md."
